{"js": "// Commit: \"CHORE: Remove some secrets\" (David Njuguna's report-generation\n// code leaked API/credential info through the SharePoint document\n// metadata it produced). The canonical-OOXML diff for this revision is\n// purely additive: three new Custom XML Parts get attached to the\n// package -- the kind of `customXml/itemN.xml` + `customXml/itemPropsN.xml`\n// pair that SharePoint/OneDrive stamps onto a file once it is associated\n// with a document-library content type:\n//\n//   1) the content-type schema (field definitions: MediaService*, Shared\n//      With*, core document properties, InfoPath \"Person\" helper schema);\n//   2) the SharePoint list-form template pointers (Display/Edit/New);\n//   3) the (empty) document-management properties part.\n//\n// `CustomXmlParts.add()` is the Word API for attaching exactly this kind\n// of package part; Word mints the paired itemPropsN.xml part (with a\n// fresh GUID), the part-relationship, and the [Content_Types].xml\n// override for every part added this way, so only the three XML payloads\n// below need to be supplied explicitly, in the same order they appear in\n// the package.\n\nconst contentTypeSchemaXml = \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<ct:contentTypeSchema xmlns:ct=\\\"http://schemas.microsoft.com/office/2006/metadata/contentType\\\" xmlns:ma=\\\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\\\" ct:_=\\\"\\\" ma:_=\\\"\\\" ma:contentTypeName=\\\"Document\\\" ma:contentTypeID=\\\"0x010100C36BF71337A24B48B564DDACE538CF46\\\" ma:contentTypeVersion=\\\"11\\\" ma:contentTypeDescription=\\\"Create a new document.\\\" ma:contentTypeScope=\\\"\\\" ma:versionID=\\\"ab11e64f46439af49bbdfd5fba0c63c9\\\">\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:p=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" xmlns:ns2=\\\"fe441e1d-e4f6-4814-912b-65f01f5d3036\\\" xmlns:ns3=\\\"3494b8ed-9bec-40f3-830b-d6b4356b7ad3\\\" targetNamespace=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" ma:root=\\\"true\\\" ma:fieldsID=\\\"89938ae16bfa147a07b2deb8c8ea1c5d\\\" ns2:_=\\\"\\\" ns3:_=\\\"\\\">\\n    <xsd:import namespace=\\\"fe441e1d-e4f6-4814-912b-65f01f5d3036\\\"/>\\n    <xsd:import namespace=\\\"3494b8ed-9bec-40f3-830b-d6b4356b7ad3\\\"/>\\n    <xsd:element name=\\\"properties\\\">\\n      <xsd:complexType>\\n        <xsd:sequence>\\n          <xsd:element name=\\\"documentManagement\\\">\\n            <xsd:complexType>\\n              <xsd:all>\\n                <xsd:element ref=\\\"ns2:MediaServiceMetadata\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceFastMetadata\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceAutoTags\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceOCR\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceGenerationTime\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceEventHashCode\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceAutoKeyPoints\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceKeyPoints\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns2:MediaServiceDateTaken\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns3:SharedWithUsers\\\" minOccurs=\\\"0\\\"/>\\n                <xsd:element ref=\\\"ns3:SharedWithDetails\\\" minOccurs=\\\"0\\\"/>\\n              </xsd:all>\\n            </xsd:complexType>\\n          </xsd:element>\\n        </xsd:sequence>\\n      </xsd:complexType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:dms=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" targetNamespace=\\\"fe441e1d-e4f6-4814-912b-65f01f5d3036\\\" elementFormDefault=\\\"qualified\\\">\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\"/>\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\"/>\\n    <xsd:element name=\\\"MediaServiceMetadata\\\" ma:index=\\\"8\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceMetadata\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceMetadata\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceFastMetadata\\\" ma:index=\\\"9\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceFastMetadata\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceFastMetadata\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceAutoTags\\\" ma:index=\\\"10\\\" nillable=\\\"true\\\" ma:displayName=\\\"Tags\\\" ma:internalName=\\\"MediaServiceAutoTags\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceOCR\\\" ma:index=\\\"11\\\" nillable=\\\"true\\\" ma:displayName=\\\"Extracted Text\\\" ma:internalName=\\\"MediaServiceOCR\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\">\\n          <xsd:maxLength value=\\\"255\\\"/>\\n        </xsd:restriction>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceGenerationTime\\\" ma:index=\\\"12\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceGenerationTime\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceGenerationTime\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceEventHashCode\\\" ma:index=\\\"13\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceEventHashCode\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceEventHashCode\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceAutoKeyPoints\\\" ma:index=\\\"14\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceAutoKeyPoints\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceAutoKeyPoints\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceKeyPoints\\\" ma:index=\\\"15\\\" nillable=\\\"true\\\" ma:displayName=\\\"KeyPoints\\\" ma:internalName=\\\"MediaServiceKeyPoints\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\">\\n          <xsd:maxLength value=\\\"255\\\"/>\\n        </xsd:restriction>\\n      </xsd:simpleType>\\n    </xsd:element>\\n    <xsd:element name=\\\"MediaServiceDateTaken\\\" ma:index=\\\"16\\\" nillable=\\\"true\\\" ma:displayName=\\\"MediaServiceDateTaken\\\" ma:hidden=\\\"true\\\" ma:internalName=\\\"MediaServiceDateTaken\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Text\\\"/>\\n      </xsd:simpleType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:dms=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" targetNamespace=\\\"3494b8ed-9bec-40f3-830b-d6b4356b7ad3\\\" elementFormDefault=\\\"qualified\\\">\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/2006/documentManagement/types\\\"/>\\n    <xsd:import namespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\"/>\\n    <xsd:element name=\\\"SharedWithUsers\\\" ma:index=\\\"17\\\" nillable=\\\"true\\\" ma:displayName=\\\"Shared With\\\" ma:internalName=\\\"SharedWithUsers\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:complexType>\\n        <xsd:complexContent>\\n          <xsd:extension base=\\\"dms:UserMulti\\\">\\n            <xsd:sequence>\\n              <xsd:element name=\\\"UserInfo\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\">\\n                <xsd:complexType>\\n                  <xsd:sequence>\\n                    <xsd:element name=\\\"DisplayName\\\" type=\\\"xsd:string\\\" minOccurs=\\\"0\\\"/>\\n                    <xsd:element name=\\\"AccountId\\\" type=\\\"dms:UserId\\\" minOccurs=\\\"0\\\" nillable=\\\"true\\\"/>\\n                    <xsd:element name=\\\"AccountType\\\" type=\\\"xsd:string\\\" minOccurs=\\\"0\\\"/>\\n                  </xsd:sequence>\\n                </xsd:complexType>\\n              </xsd:element>\\n            </xsd:sequence>\\n          </xsd:extension>\\n        </xsd:complexContent>\\n      </xsd:complexType>\\n    </xsd:element>\\n    <xsd:element name=\\\"SharedWithDetails\\\" ma:index=\\\"18\\\" nillable=\\\"true\\\" ma:displayName=\\\"Shared With Details\\\" ma:internalName=\\\"SharedWithDetails\\\" ma:readOnly=\\\"true\\\">\\n      <xsd:simpleType>\\n        <xsd:restriction base=\\\"dms:Note\\\">\\n          <xsd:maxLength value=\\\"255\\\"/>\\n        </xsd:restriction>\\n      </xsd:simpleType>\\n    </xsd:element>\\n  </xsd:schema>\\n  <xsd:schema xmlns=\\\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\\\" xmlns:xsd=\\\"http://www.w3.org/2001/XMLSchema\\\" xmlns:xsi=\\\"http://www.w3.org/2001/XMLSchema-instance\\\" xmlns:dc=\\\"http://purl.org/dc/elements/1.1/\\\" xmlns:dcterms=\\\"http://purl.org/dc/terms/\\\" xmlns:odoc=\\\"http://schemas.microsoft.com/internal/obd\\\" targetNamespace=\\\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\\\" elementFormDefault=\\\"qualified\\\" attributeFormDefault=\\\"unqualified\\\" blockDefault=\\\"#all\\\">\\n    <xsd:import namespace=\\\"http://purl.org/dc/elements/1.1/\\\" schemaLocation=\\\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\\\"/>\\n    <xsd:import namespace=\\\"http://purl.org/dc/terms/\\\" schemaLocation=\\\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\\\"/>\\n    <xsd:element name=\\\"coreProperties\\\" type=\\\"CT_coreProperties\\\"/>\\n    <xsd:complexType name=\\\"CT_coreProperties\\\">\\n      <xsd:all>\\n        <xsd:element ref=\\\"dc:creator\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dcterms:created\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dc:identifier\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"contentType\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\" ma:index=\\\"0\\\" ma:displayName=\\\"Content Type\\\"/>\\n        <xsd:element ref=\\\"dc:title\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" ma:index=\\\"4\\\" ma:displayName=\\\"Title\\\"/>\\n        <xsd:element ref=\\\"dc:subject\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element ref=\\\"dc:description\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"keywords\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element ref=\\\"dc:language\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"category\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element name=\\\"version\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element name=\\\"revision\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\">\\n          <xsd:annotation>\\n            <xsd:documentation>\\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\\n                    </xsd:documentation>\\n          </xsd:annotation>\\n        </xsd:element>\\n        <xsd:element name=\\\"lastModifiedBy\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n        <xsd:element ref=\\\"dcterms:modified\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\"/>\\n        <xsd:element name=\\\"contentStatus\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"1\\\" type=\\\"xsd:string\\\"/>\\n      </xsd:all>\\n    </xsd:complexType>\\n  </xsd:schema>\\n  <xs:schema xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" xmlns:xs=\\\"http://www.w3.org/2001/XMLSchema\\\" targetNamespace=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\" elementFormDefault=\\\"qualified\\\" attributeFormDefault=\\\"unqualified\\\">\\n    <xs:element name=\\\"Person\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:DisplayName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:AccountId\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:AccountType\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"DisplayName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"AccountId\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"AccountType\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"BDCAssociatedEntity\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:BDCEntity\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\"/>\\n        </xs:sequence>\\n        <xs:attribute ref=\\\"pc:EntityNamespace\\\"/>\\n        <xs:attribute ref=\\\"pc:EntityName\\\"/>\\n        <xs:attribute ref=\\\"pc:SystemInstanceName\\\"/>\\n        <xs:attribute ref=\\\"pc:AssociationName\\\"/>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:attribute name=\\\"EntityNamespace\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"EntityName\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"SystemInstanceName\\\" type=\\\"xs:string\\\"/>\\n    <xs:attribute name=\\\"AssociationName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"BDCEntity\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:EntityDisplayName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityInstanceReference\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId1\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId2\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId3\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId4\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:EntityId5\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"EntityDisplayName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityInstanceReference\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId1\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId2\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId3\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId4\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"EntityId5\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"Terms\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:TermInfo\\\" minOccurs=\\\"0\\\" maxOccurs=\\\"unbounded\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"TermInfo\\\">\\n      <xs:complexType>\\n        <xs:sequence>\\n          <xs:element ref=\\\"pc:TermName\\\" minOccurs=\\\"0\\\"/>\\n          <xs:element ref=\\\"pc:TermId\\\" minOccurs=\\\"0\\\"/>\\n        </xs:sequence>\\n      </xs:complexType>\\n    </xs:element>\\n    <xs:element name=\\\"TermName\\\" type=\\\"xs:string\\\"/>\\n    <xs:element name=\\\"TermId\\\" type=\\\"xs:string\\\"/>\\n  </xs:schema>\\n</ct:contentTypeSchema>\\n\";\n\nconst formTemplatesXml = \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<?mso-contentType ?>\\n<FormTemplates xmlns=\\\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\\\">\\n  <Display>DocumentLibraryForm</Display>\\n  <Edit>DocumentLibraryForm</Edit>\\n  <New>DocumentLibraryForm</New>\\n</FormTemplates>\\n\";\n\nconst documentManagementPropsXml = \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\"?>\\n<p:properties xmlns:p=\\\"http://schemas.microsoft.com/office/2006/metadata/properties\\\" xmlns:xsi=\\\"http://www.w3.org/2001/XMLSchema-instance\\\" xmlns:pc=\\\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\\\">\\n  <documentManagement/>\\n</p:properties>\\n\";\n\nconst customXmlParts = context.document.customXmlParts;\ncustomXmlParts.add(contentTypeSchemaXml);\ncustomXmlParts.add(formTemplatesXml);\ncustomXmlParts.add(documentManagementPropsXml);\n\nawait context.sync();\n", "ps1": "# Commit: \"CHORE: Remove some secrets\" (David Njuguna's report-generation\n# code leaked API/credential info through the SharePoint document\n# metadata it produced). The canonical-OOXML diff for this revision is\n# purely additive: three new Custom XML Parts get attached to the\n# package -- the kind of customXml/itemN.xml + customXml/itemPropsN.xml\n# pair that SharePoint/OneDrive stamps onto a file once it is associated\n# with a document-library content type:\n#\n#   1) the content-type schema (field definitions: MediaService*, Shared\n#      With*, core document properties, InfoPath \"Person\" helper schema);\n#   2) the SharePoint list-form template pointers (Display/Edit/New);\n#   3) the (empty) document-management properties part.\n#\n# CustomXMLParts.Add is the Word COM API for attaching exactly this kind\n# of package part; Word mints the paired itemPropsN.xml part (with a\n# fresh GUID), the part-relationship, and the [Content_Types].xml\n# override for every part added this way, so only the three XML payloads\n# below need to be supplied explicitly, in the same order they appear in\n# the package.\n\n$d = $word.ActiveDocument\n\n$contentTypeSchemaXml = '<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<ct:contentTypeSchema xmlns:ct=\"http://schemas.microsoft.com/office/2006/metadata/contentType\" xmlns:ma=\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\" ct:_=\"\" ma:_=\"\" ma:contentTypeName=\"Document\" ma:contentTypeID=\"0x010100C36BF71337A24B48B564DDACE538CF46\" ma:contentTypeVersion=\"11\" ma:contentTypeDescription=\"Create a new document.\" ma:contentTypeScope=\"\" ma:versionID=\"ab11e64f46439af49bbdfd5fba0c63c9\">\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:ns2=\"fe441e1d-e4f6-4814-912b-65f01f5d3036\" xmlns:ns3=\"3494b8ed-9bec-40f3-830b-d6b4356b7ad3\" targetNamespace=\"http://schemas.microsoft.com/office/2006/metadata/properties\" ma:root=\"true\" ma:fieldsID=\"89938ae16bfa147a07b2deb8c8ea1c5d\" ns2:_=\"\" ns3:_=\"\">\n    <xsd:import namespace=\"fe441e1d-e4f6-4814-912b-65f01f5d3036\"/>\n    <xsd:import namespace=\"3494b8ed-9bec-40f3-830b-d6b4356b7ad3\"/>\n    <xsd:element name=\"properties\">\n      <xsd:complexType>\n        <xsd:sequence>\n          <xsd:element name=\"documentManagement\">\n            <xsd:complexType>\n              <xsd:all>\n                <xsd:element ref=\"ns2:MediaServiceMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceFastMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoTags\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceOCR\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceGenerationTime\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceEventHashCode\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceDateTaken\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithUsers\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithDetails\" minOccurs=\"0\"/>\n              </xsd:all>\n            </xsd:complexType>\n          </xsd:element>\n        </xsd:sequence>\n      </xsd:complexType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"fe441e1d-e4f6-4814-912b-65f01f5d3036\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"MediaServiceMetadata\" ma:index=\"8\" nillable=\"true\" ma:displayName=\"MediaServiceMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceFastMetadata\" ma:index=\"9\" nillable=\"true\" ma:displayName=\"MediaServiceFastMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceFastMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoTags\" ma:index=\"10\" nillable=\"true\" ma:displayName=\"Tags\" ma:internalName=\"MediaServiceAutoTags\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceOCR\" ma:index=\"11\" nillable=\"true\" ma:displayName=\"Extracted Text\" ma:internalName=\"MediaServiceOCR\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceGenerationTime\" ma:index=\"12\" nillable=\"true\" ma:displayName=\"MediaServiceGenerationTime\" ma:hidden=\"true\" ma:internalName=\"MediaServiceGenerationTime\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceEventHashCode\" ma:index=\"13\" nillable=\"true\" ma:displayName=\"MediaServiceEventHashCode\" ma:hidden=\"true\" ma:internalName=\"MediaServiceEventHashCode\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoKeyPoints\" ma:index=\"14\" nillable=\"true\" ma:displayName=\"MediaServiceAutoKeyPoints\" ma:hidden=\"true\" ma:internalName=\"MediaServiceAutoKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceKeyPoints\" ma:index=\"15\" nillable=\"true\" ma:displayName=\"KeyPoints\" ma:internalName=\"MediaServiceKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceDateTaken\" ma:index=\"16\" nillable=\"true\" ma:displayName=\"MediaServiceDateTaken\" ma:hidden=\"true\" ma:internalName=\"MediaServiceDateTaken\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"3494b8ed-9bec-40f3-830b-d6b4356b7ad3\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"SharedWithUsers\" ma:index=\"17\" nillable=\"true\" ma:displayName=\"Shared With\" ma:internalName=\"SharedWithUsers\" ma:readOnly=\"true\">\n      <xsd:complexType>\n        <xsd:complexContent>\n          <xsd:extension base=\"dms:UserMulti\">\n            <xsd:sequence>\n              <xsd:element name=\"UserInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\">\n                <xsd:complexType>\n                  <xsd:sequence>\n                    <xsd:element name=\"DisplayName\" type=\"xsd:string\" minOccurs=\"0\"/>\n                    <xsd:element name=\"AccountId\" type=\"dms:UserId\" minOccurs=\"0\" nillable=\"true\"/>\n                    <xsd:element name=\"AccountType\" type=\"xsd:string\" minOccurs=\"0\"/>\n                  </xsd:sequence>\n                </xsd:complexType>\n              </xsd:element>\n            </xsd:sequence>\n          </xsd:extension>\n        </xsd:complexContent>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"SharedWithDetails\" ma:index=\"18\" nillable=\"true\" ma:displayName=\"Shared With Details\" ma:internalName=\"SharedWithDetails\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:dc=\"http://purl.org/dc/elements/1.1/\" xmlns:dcterms=\"http://purl.org/dc/terms/\" xmlns:odoc=\"http://schemas.microsoft.com/internal/obd\" targetNamespace=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\" blockDefault=\"#all\">\n    <xsd:import namespace=\"http://purl.org/dc/elements/1.1/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\"/>\n    <xsd:import namespace=\"http://purl.org/dc/terms/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\"/>\n    <xsd:element name=\"coreProperties\" type=\"CT_coreProperties\"/>\n    <xsd:complexType name=\"CT_coreProperties\">\n      <xsd:all>\n        <xsd:element ref=\"dc:creator\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dcterms:created\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:identifier\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentType\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\" ma:index=\"0\" ma:displayName=\"Content Type\"/>\n        <xsd:element ref=\"dc:title\" minOccurs=\"0\" maxOccurs=\"1\" ma:index=\"4\" ma:displayName=\"Title\"/>\n        <xsd:element ref=\"dc:subject\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:description\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"keywords\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dc:language\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"category\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"version\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"revision\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\">\n          <xsd:annotation>\n            <xsd:documentation>\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\n                    </xsd:documentation>\n          </xsd:annotation>\n        </xsd:element>\n        <xsd:element name=\"lastModifiedBy\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dcterms:modified\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentStatus\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n      </xsd:all>\n    </xsd:complexType>\n  </xsd:schema>\n  <xs:schema xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" targetNamespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\">\n    <xs:element name=\"Person\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:DisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountId\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountType\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"DisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"AccountId\" type=\"xs:string\"/>\n    <xs:element name=\"AccountType\" type=\"xs:string\"/>\n    <xs:element name=\"BDCAssociatedEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:BDCEntity\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n        <xs:attribute ref=\"pc:EntityNamespace\"/>\n        <xs:attribute ref=\"pc:EntityName\"/>\n        <xs:attribute ref=\"pc:SystemInstanceName\"/>\n        <xs:attribute ref=\"pc:AssociationName\"/>\n      </xs:complexType>\n    </xs:element>\n    <xs:attribute name=\"EntityNamespace\" type=\"xs:string\"/>\n    <xs:attribute name=\"EntityName\" type=\"xs:string\"/>\n    <xs:attribute name=\"SystemInstanceName\" type=\"xs:string\"/>\n    <xs:attribute name=\"AssociationName\" type=\"xs:string\"/>\n    <xs:element name=\"BDCEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:EntityDisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityInstanceReference\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId1\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId2\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId3\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId4\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId5\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"EntityDisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"EntityInstanceReference\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId1\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId2\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId3\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId4\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId5\" type=\"xs:string\"/>\n    <xs:element name=\"Terms\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermInfo\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:TermId\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermName\" type=\"xs:string\"/>\n    <xs:element name=\"TermId\" type=\"xs:string\"/>\n  </xs:schema>\n</ct:contentTypeSchema>\n'\n\n$formTemplatesXml = '<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<?mso-contentType ?>\n<FormTemplates xmlns=\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\">\n  <Display>DocumentLibraryForm</Display>\n  <Edit>DocumentLibraryForm</Edit>\n  <New>DocumentLibraryForm</New>\n</FormTemplates>\n'\n\n$documentManagementPropsXml = '<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<p:properties xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\">\n  <documentManagement/>\n</p:properties>\n'\n\n$d.CustomXMLParts.Add($contentTypeSchemaXml) | Out-Null\n$d.CustomXMLParts.Add($formTemplatesXml) | Out-Null\n$d.CustomXMLParts.Add($documentManagementPropsXml) | Out-Null\n"}
